{"js": "// Products.docx clean-up after the HTML -> .docx conversion:\n//\n//  1. The converter emitted THREE consecutive centered, empty, 48pt\n//     paragraphs right after the \"Products\" title text (the title run\n//     plus two blank duplicate paragraphs). Only one blank centered\n//     paragraph should remain, so drop the extra (first) one.\n//\n//  2. Two bulleted paragraphs (\"for apps on the App Store, enter: ...\"\n//     and \"for apps on the google play store, enter: ...\") carry a\n//     bogus negative hanging indent left over from the conversion\n//     (left=720, first-line=-720). They should keep only the\n//     left = 720 (36pt) indent.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"alignment,text,leftIndent\");\nawait context.sync();\n\n// --- 1. Remove the extra blank centered paragraph ---\n// Find the first paragraph that is centered and has no visible text;\n// that is the stray duplicate sitting right after the \"Products\"\n// heading paragraph (the heading itself contains the word \"Products\",\n// so it is skipped).\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.alignment === Word.Alignment.centered && p.text.trim() === \"\") {\n    p.delete();\n    await context.sync();\n    break;\n  }\n}\n\n// --- 2. Fix the two hanging-indent bullet paragraphs ---\n// Re-fetch the (now shifted) paragraph collection and re-assign\n// leftIndent to its own value on every paragraph that carries the\n// 36pt (720 twip) left indent; that forces Word to rewrite the\n// paragraph's indentation from the left=36pt model, dropping the\n// stray first-line indent the converter left behind.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nparagraphs2.load(\"leftIndent\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const p = paragraphs2.items[i];\n  if (p.leftIndent === 36) {\n    p.leftIndent = 36;\n  }\n}\nawait context.sync();\n", "ps1": "# Products.docx clean-up after the HTML -> .docx conversion:\n#\n#  1. The converter emitted THREE consecutive centered, empty, 48pt\n#     paragraphs right after the \"Products\" title text (the title run\n#     plus two blank duplicate paragraphs). Only one blank centered\n#     paragraph should remain, so drop the extra (first) one.\n#\n#  2. Two bulleted paragraphs (\"for apps on the App Store, enter: ...\"\n#     and \"for apps on the google play store, enter: ...\") carry a\n#     bogus negative hanging indent left over from the conversion\n#     (left=720, first-line=-720). They should keep only the\n#     w:left=\"720\" indent.\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the extra blank centered paragraph ---\n# Walk the document looking for the first paragraph that is centered\n# and has no visible text; that is the stray duplicate sitting right\n# after the \"Products\" heading paragraph (the heading itself contains\n# the word \"Products\", so it is skipped).\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Format.Alignment -eq 1 -and $p.Range.Text.Trim() -eq \"\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 2. Fix the two hanging-indent bullet paragraphs ---\n# Re-assigning LeftIndent to its own value forces Word to rewrite the\n# paragraph's <w:ind> from the (left=720) model, which drops the stray\n# first-line attribute the converter left behind.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Format.LeftIndent -eq 36) {\n        $p.Format.LeftIndent = $p.Format.LeftIndent\n    }\n}\n"}
